$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("us")

function Set-TextValue($addr, $val) {
    # Cells in column M store numeric-looking values as literal text
    # (dividend yields exported verbatim from a data source). Force the
    # "@" text format before assigning so Excel keeps it as a string
    # instead of silently parsing it into a number, then strip the
    # temporary formatting so no stray style is left behind on the cell.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

# Dividend yield updates (column M) - stored as text
Set-TextValue "M2" "0.131"
Set-TextValue "M3" "0.1348"
Set-TextValue "M8" "0.1349"
Set-TextValue "M11" "0.2373"
Set-TextValue "M12" "0.12020001"
Set-TextValue "M16" "0.1224"
Set-TextValue "M17" "0.119899996"
Set-TextValue "M22" "0.1741"
Set-TextValue "M27" "0.1481"
Set-TextValue "M28" "0.1617"
Set-TextValue "M37" "0.1844"
Set-TextValue "M54" "0.1187"
Set-TextValue "M56" "0.0771"
Set-TextValue "M57" "0.1291"

# Current price update (column G) - numeric
$ws.Range("G24").Value = 8.25

# Indicator updates (column C)
$ws.Range("C30").Value = "MONITOR"
$ws.Range("C33").Value = "MONITOR"
$ws.Range("C38").Value = "MONITOR"
$ws.Range("C64").Value = "MONITOR"

# Earnings date / estimate dict-like strings (column R)
$ws.Range("R7").Value = "{'earningsDate': ['2023-08-03'], 'earningsAverage': 0.32, 'earningsLow': 0.13, 'earningsHigh': 0.43, 'revenueAverage': 119161000000, 'revenueLow': 115113000000, 'revenueHigh': 120834000000}"
$ws.Range("R13").Value = "{'earningsDate': ['2023-10-24'], 'earningsAverage': -0.04, 'earningsLow': -0.04, 'earningsHigh': -0.04, 'revenueAverage': 130970000, 'revenueLow': 128000000, 'revenueHigh': 133000000}"
$ws.Range("R38").Value = "{'earningsDate': [], 'earningsAverage': 0.0, 'earningsLow': 0.0, 'earningsHigh': 0.0, 'revenueAverage': 202350000, 'revenueLow': 199220000, 'revenueHigh': 205480000}"
$ws.Range("R62").Value = "{'earningsDate': ['2023-10-17', '2023-10-23'], 'earningsAverage': 0.73, 'earningsLow': 0.59, 'earningsHigh': 0.92, 'revenueAverage': 22401900000, 'revenueLow': 18372700000, 'revenueHigh': 24730200000}"
